# Insert two new rows of weekly price data for "Coliflor" right after the
# current row 490, pushing all subsequent rows (old 491..519) down to 493..521.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("491:492").Insert()

# New row 491
$ws.Cells.Item(491, 1).Value = 11
$ws.Cells.Item(491, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(491, 3).Value = "Bíobío"
$ws.Cells.Item(491, 4).Value = "12/7/2023"
$ws.Cells.Item(491, 5).Value = 8
$ws.Cells.Item(491, 6).Value = 100112008
$ws.Cells.Item(491, 7).Value = "Coliflor"
$ws.Cells.Item(491, 8).Value = "Sin especificar"
$ws.Cells.Item(491, 9).Value = "Primera"
$ws.Cells.Item(491, 10).Value = 1000
$ws.Cells.Item(491, 11).Value = 1000
$ws.Cells.Item(491, 12).Value = 1000
$ws.Cells.Item(491, 13).Value = 1000
$ws.Cells.Item(491, 14).Value = "$/unidad"
$ws.Cells.Item(491, 15).Value = "Región Metropolitana"
$ws.Cells.Item(491, 16).Value = 1000
$ws.Cells.Item(491, 17).Value = 1
$ws.Cells.Item(491, 18).Value = "Hortaliza"

# New row 492
$ws.Cells.Item(492, 1).Value = 11
$ws.Cells.Item(492, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(492, 3).Value = "Bíobío"
$ws.Cells.Item(492, 4).Value = "12/7/2023"
$ws.Cells.Item(492, 5).Value = 8
$ws.Cells.Item(492, 6).Value = 100112008
$ws.Cells.Item(492, 7).Value = "Coliflor"
$ws.Cells.Item(492, 8).Value = "Sin especificar"
$ws.Cells.Item(492, 9).Value = "Segunda"
$ws.Cells.Item(492, 10).Value = 1000
$ws.Cells.Item(492, 11).Value = 700
$ws.Cells.Item(492, 12).Value = 700
$ws.Cells.Item(492, 13).Value = 700
$ws.Cells.Item(492, 14).Value = "$/unidad"
$ws.Cells.Item(492, 15).Value = "Región Metropolitana"
$ws.Cells.Item(492, 16).Value = 700
$ws.Cells.Item(492, 17).Value = 1
$ws.Cells.Item(492, 18).Value = "Hortaliza"
